$d = $word.ActiveDocument

# --- Row 1 of the table: cell-by-cell rewrite (one cell removed, one cell added) ---
$t = $d.Tables.Item(1)
$cell = $t.Cell(1,1)
$r = $cell.Range
$d.Range($r.Start, $r.End).Text = "59-50=9"

$cell = $t.Cell(1,2)
$r = $cell.Range
$d.Range($r.Start, $r.End).Text = "33-22=11"

$cell = $t.Cell(1,3)
$r = $cell.Range
$d.Range($r.Start, $r.End).Text = "68+30=98"

$cell = $t.Cell(1,4)
$r = $cell.Range
$d.Range($r.Start, $r.End).Text = "62-26=36"

$cell = $t.Cell(1,5)
$r = $cell.Range
$d.Range($r.Start, $r.End).Text = "85-37=48"

# --- Remaining cells: in-place Find & Replace (unique text per cell) ---
$d.Content.Find.Execute("87-78=9", $true, $false, $false, $false, $false, $true, 1, $false, "19+48=67", 2) | Out-Null
$d.Content.Find.Execute("67+28=95", $true, $false, $false, $false, $false, $true, 1, $false, "34+23=57", 2) | Out-Null
$d.Content.Find.Execute("63-22=41", $true, $false, $false, $false, $false, $true, 1, $false, "76-50=26", 2) | Out-Null
$d.Content.Find.Execute("97-73=24", $true, $false, $false, $false, $false, $true, 1, $false, "79-62=17", 2) | Out-Null
$d.Content.Find.Execute("49-41=8", $true, $false, $false, $false, $false, $true, 1, $false, "51-34=17", 2) | Out-Null
$d.Content.Find.Execute("42+15=57", $true, $false, $false, $false, $false, $true, 1, $false, "31+23=54", 2) | Out-Null
$d.Content.Find.Execute("75-34=41", $true, $false, $false, $false, $false, $true, 1, $false, "22+42=64", 2) | Out-Null
$d.Content.Find.Execute("26+51=77", $true, $false, $false, $false, $false, $true, 1, $false, "15-2=13", 2) | Out-Null
$d.Content.Find.Execute("81-2=79", $true, $false, $false, $false, $false, $true, 1, $false, "84-57=27", 2) | Out-Null
$d.Content.Find.Execute("28+5=33", $true, $false, $false, $false, $false, $true, 1, $false, "9+23=32", 2) | Out-Null
$d.Content.Find.Execute("0+5=5", $true, $false, $false, $false, $false, $true, 1, $false, "10+18=28", 2) | Out-Null
$d.Content.Find.Execute("44-31=13", $true, $false, $false, $false, $false, $true, 1, $false, "54-16=38", 2) | Out-Null
$d.Content.Find.Execute("14+67=81", $true, $false, $false, $false, $false, $true, 1, $false, "43-35=8", 2) | Out-Null
$d.Content.Find.Execute("52+37=89", $true, $false, $false, $false, $false, $true, 1, $false, "60+18=78", 2) | Out-Null
$d.Content.Find.Execute("42+23=65", $true, $false, $false, $false, $false, $true, 1, $false, "80-7=73", 2) | Out-Null
$d.Content.Find.Execute("50-45=5", $true, $false, $false, $false, $false, $true, 1, $false, "75-23=52", 2) | Out-Null
$d.Content.Find.Execute("10-1=9", $true, $false, $false, $false, $false, $true, 1, $false, "98-39=59", 2) | Out-Null
$d.Content.Find.Execute("61-56=5", $true, $false, $false, $false, $false, $true, 1, $false, "87-38=49", 2) | Out-Null
$d.Content.Find.Execute("27+41=68", $true, $false, $false, $false, $false, $true, 1, $false, "24-22=2", 2) | Out-Null
$d.Content.Find.Execute("85-18=67", $true, $false, $false, $false, $false, $true, 1, $false, "0+49=49", 2) | Out-Null
$d.Content.Find.Execute("80+1=81", $true, $false, $false, $false, $false, $true, 1, $false, "2+92=94", 2) | Out-Null
$d.Content.Find.Execute("29+34=63", $true, $false, $false, $false, $false, $true, 1, $false, "47-34=13", 2) | Out-Null
$d.Content.Find.Execute("36+49=85", $true, $false, $false, $false, $false, $true, 1, $false, "14+29=43", 2) | Out-Null
$d.Content.Find.Execute("85-73=12", $true, $false, $false, $false, $false, $true, 1, $false, "44+35=79", 2) | Out-Null
$d.Content.Find.Execute("5+21=26", $true, $false, $false, $false, $false, $true, 1, $false, "71-47=24", 2) | Out-Null
$d.Content.Find.Execute("97-83=14", $true, $false, $false, $false, $false, $true, 1, $false, "83-65=18", 2) | Out-Null
$d.Content.Find.Execute("71-11=60", $true, $false, $false, $false, $false, $true, 1, $false, "39+30=69", 2) | Out-Null
$d.Content.Find.Execute("5+29=34", $true, $false, $false, $false, $false, $true, 1, $false, "94-90=4", 2) | Out-Null
$d.Content.Find.Execute("24-10=14", $true, $false, $false, $false, $false, $true, 1, $false, "29+23=52", 2) | Out-Null
$d.Content.Find.Execute("74+16=90", $true, $false, $false, $false, $false, $true, 1, $false, "55+36=91", 2) | Out-Null
$d.Content.Find.Execute("1+3=4", $true, $false, $false, $false, $false, $true, 1, $false, "57-53=4", 2) | Out-Null
$d.Content.Find.Execute("64-29=35", $true, $false, $false, $false, $false, $true, 1, $false, "63+8=71", 2) | Out-Null
$d.Content.Find.Execute("44-7=37", $true, $false, $false, $false, $false, $true, 1, $false, "82-9=73", 2) | Out-Null
$d.Content.Find.Execute("82-8=74", $true, $false, $false, $false, $false, $true, 1, $false, "41+1=42", 2) | Out-Null
$d.Content.Find.Execute("28+8=36", $true, $false, $false, $false, $false, $true, 1, $false, "37-27=10", 2) | Out-Null
$d.Content.Find.Execute("28+59=87", $true, $false, $false, $false, $false, $true, 1, $false, "1+62=63", 2) | Out-Null
$d.Content.Find.Execute("44+0=44", $true, $false, $false, $false, $false, $true, 1, $false, "65-48=17", 2) | Out-Null
$d.Content.Find.Execute("73+16=89", $true, $false, $false, $false, $false, $true, 1, $false, "9+8=17", 2) | Out-Null
$d.Content.Find.Execute("39+41=80", $true, $false, $false, $false, $false, $true, 1, $false, "7+2=9", 2) | Out-Null
$d.Content.Find.Execute("91-26=65", $true, $false, $false, $false, $false, $true, 1, $false, "58+36=94", 2) | Out-Null
$d.Content.Find.Execute("4+57=61", $true, $false, $false, $false, $false, $true, 1, $false, "30+40=70", 2) | Out-Null
$d.Content.Find.Execute("0+31=31", $true, $false, $false, $false, $false, $true, 1, $false, "20+67=87", 2) | Out-Null
$d.Content.Find.Execute("63-10=53", $true, $false, $false, $false, $false, $true, 1, $false, "6-0=6", 2) | Out-Null
$d.Content.Find.Execute("98-85=13", $true, $false, $false, $false, $false, $true, 1, $false, "68+24=92", 2) | Out-Null
$d.Content.Find.Execute("20+66=86", $true, $false, $false, $false, $false, $true, 1, $false, "76-45=31", 2) | Out-Null
$d.Content.Find.Execute("21-9=12", $true, $false, $false, $false, $false, $true, 1, $false, "29-5=24", 2) | Out-Null
$d.Content.Find.Execute("11+8=19", $true, $false, $false, $false, $false, $true, 1, $false, "49-18=31", 2) | Out-Null
$d.Content.Find.Execute("55+3=58", $true, $false, $false, $false, $false, $true, 1, $false, "76-70=6", 2) | Out-Null
$d.Content.Find.Execute("57+42=99", $true, $false, $false, $false, $false, $true, 1, $false, "89-47=42", 2) | Out-Null
$d.Content.Find.Execute("65+7=72", $true, $false, $false, $false, $false, $true, 1, $false, "99-40=59", 2) | Out-Null
$d.Content.Find.Execute("7+73=80", $true, $false, $false, $false, $false, $true, 1, $false, "23+41=64", 2) | Out-Null
$d.Content.Find.Execute("34-24=10", $true, $false, $false, $false, $false, $true, 1, $false, "37-17=20", 2) | Out-Null
$d.Content.Find.Execute("27-16=11", $true, $false, $false, $false, $false, $true, 1, $false, "68-8=60", 2) | Out-Null
$d.Content.Find.Execute("40-7=33", $true, $false, $false, $false, $false, $true, 1, $false, "82+15=97", 2) | Out-Null
$d.Content.Find.Execute("56+37=93", $true, $false, $false, $false, $false, $true, 1, $false, "22+61=83", 2) | Out-Null
$d.Content.Find.Execute("88+10=98", $true, $false, $false, $false, $false, $true, 1, $false, "55+19=74", 2) | Out-Null
$d.Content.Find.Execute("61-34=27", $true, $false, $false, $false, $false, $true, 1, $false, "35-29=6", 2) | Out-Null
$d.Content.Find.Execute("29+68=97", $true, $false, $false, $false, $false, $true, 1, $false, "62-17=45", 2) | Out-Null
$d.Content.Find.Execute("15+57=72", $true, $false, $false, $false, $false, $true, 1, $false, "39+25=64", 2) | Out-Null
$d.Content.Find.Execute("61+33=94", $true, $false, $false, $false, $false, $true, 1, $false, "5+48=53", 2) | Out-Null
$d.Content.Find.Execute("84+3=87", $true, $false, $false, $false, $false, $true, 1, $false, "71+15=86", 2) | Out-Null
$d.Content.Find.Execute("93-81=12", $true, $false, $false, $false, $false, $true, 1, $false, "27-18=9", 2) | Out-Null
$d.Content.Find.Execute("4+81=85", $true, $false, $false, $false, $false, $true, 1, $false, "86+11=97", 2) | Out-Null
$d.Content.Find.Execute("79-69=10", $true, $false, $false, $false, $false, $true, 1, $false, "96-62=34", 2) | Out-Null
$d.Content.Find.Execute("27-24=3", $true, $false, $false, $false, $false, $true, 1, $false, "9+55=64", 2) | Out-Null
$d.Content.Find.Execute("56-11=45", $true, $false, $false, $false, $false, $true, 1, $false, "41-1=40", 2) | Out-Null
$d.Content.Find.Execute("10+12=22", $true, $false, $false, $false, $false, $true, 1, $false, "71+8=79", 2) | Out-Null
$d.Content.Find.Execute("11+54=65", $true, $false, $false, $false, $false, $true, 1, $false, "56-17=39", 2) | Out-Null
$d.Content.Find.Execute("70-47=23", $true, $false, $false, $false, $false, $true, 1, $false, "59-25=34", 2) | Out-Null
$d.Content.Find.Execute("77-53=24", $true, $false, $false, $false, $false, $true, 1, $false, "69-65=4", 2) | Out-Null
$d.Content.Find.Execute("87-45=42", $true, $false, $false, $false, $false, $true, 1, $false, "13+21=34", 2) | Out-Null
$d.Content.Find.Execute("3+81=84", $true, $false, $false, $false, $false, $true, 1, $false, "27-19=8", 2) | Out-Null
$d.Content.Find.Execute("53+12=65", $true, $false, $false, $false, $false, $true, 1, $false, "87-87=0", 2) | Out-Null
$d.Content.Find.Execute("34+52=86", $true, $false, $false, $false, $false, $true, 1, $false, "33-15=18", 2) | Out-Null
$d.Content.Find.Execute("35+13=48", $true, $false, $false, $false, $false, $true, 1, $false, "86-16=70", 2) | Out-Null
$d.Content.Find.Execute("84-79=5", $true, $false, $false, $false, $false, $true, 1, $false, "50+32=82", 2) | Out-Null
$d.Content.Find.Execute("87-86=1", $true, $false, $false, $false, $false, $true, 1, $false, "13+56=69", 2) | Out-Null
$d.Content.Find.Execute("63-62=1", $true, $false, $false, $false, $false, $true, 1, $false, "73+26=99", 2) | Out-Null
$d.Content.Find.Execute("1+87=88", $true, $false, $false, $false, $false, $true, 1, $false, "71-47=24", 2) | Out-Null
$d.Content.Find.Execute("93-60=33", $true, $false, $false, $false, $false, $true, 1, $false, "81-44=37", 2) | Out-Null
$d.Content.Find.Execute("34-17=17", $true, $false, $false, $false, $false, $true, 1, $false, "61+6=67", 2) | Out-Null
$d.Content.Find.Execute("22+48=70", $true, $false, $false, $false, $false, $true, 1, $false, "74+5=79", 2) | Out-Null
$d.Content.Find.Execute("76-41=35", $true, $false, $false, $false, $false, $true, 1, $false, "60-36=24", 2) | Out-Null
$d.Content.Find.Execute("96-41=55", $true, $false, $false, $false, $false, $true, 1, $false, "35-11=24", 2) | Out-Null
$d.Content.Find.Execute("59-28=31", $true, $false, $false, $false, $false, $true, 1, $false, "33+59=92", 2) | Out-Null
$d.Content.Find.Execute("74+1=75", $true, $false, $false, $false, $false, $true, 1, $false, "56-18=38", 2) | Out-Null
$d.Content.Find.Execute("72-59=13", $true, $false, $false, $false, $false, $true, 1, $false, "71+25=96", 2) | Out-Null
$d.Content.Find.Execute("29+59=88", $true, $false, $false, $false, $false, $true, 1, $false, "26+39=65", 2) | Out-Null
$d.Content.Find.Execute("67-43=24", $true, $false, $false, $false, $false, $true, 1, $false, "92+1=93", 2) | Out-Null
$d.Content.Find.Execute("78-68=10", $true, $false, $false, $false, $false, $true, 1, $false, "82-74=8", 2) | Out-Null
$d.Content.Find.Execute("7+37=44", $true, $false, $false, $false, $false, $true, 1, $false, "18+35=53", 2) | Out-Null
$d.Content.Find.Execute("51-50=1", $true, $false, $false, $false, $false, $true, 1, $false, "17+7=24", 2) | Out-Null
$d.Content.Find.Execute("24+14=38", $true, $false, $false, $false, $false, $true, 1, $false, "67+14=81", 2) | Out-Null
$d.Content.Find.Execute("87-70=17", $true, $false, $false, $false, $false, $true, 1, $false, "18+75=93", 2) | Out-Null
$d.Content.Find.Execute("86-57=29", $true, $false, $false, $false, $false, $true, 1, $false, "3+10=13", 2) | Out-Null

Write-Output "done"